$d = $word.ActiveDocument

# 1. Append the new sentence to the existing run's text. A temporary
#    trailing marker character is appended too, so that the insertion
#    point used below for the bookmark is not the very last character
#    of the paragraph (placing a bookmark exactly before a paragraph
#    mark is unreliable), then the marker is stripped back out.
$d.Content.Find.Execute(
    "多云转小雨，今天学习了分支管理，创建了一个dev分支。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "多云转小雨，今天学习了分支管理，创建了一个dev分支。使用GIT创建分支简单又快速。@",
    2)

# 2. Find the paragraph that now holds this sentence.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*使用GIT创建分支简单又快速。@*") {
        $targetPara = $p
        break
    }
}

# 3. Move the insertion point to just before the temporary marker
#    (i.e. right after the text, still before the paragraph mark) and
#    (re)create the "_GoBack" bookmark there. Adding a bookmark with an
#    already-existing name moves it, so the stray one further down in
#    the document collapses into this new location.
$markerPos = $targetPara.Range.End - 2
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 4. Remove the temporary marker character.
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()
